$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.706.47"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "2.891.46"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.38"
$ws.Range("E5").Value = "  -4.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.86"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("D9").Value = "2.890.83"
$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.62"
$ws.Range("E14").Value = "  -2.39%  "

$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").Value = "3.371.12"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").Value = "61.655.30"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "2.900.29"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").Value = "  -2.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.97"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.99"
$ws.Range("E21").Value = "  -2.86%  "

$ws.Range("E22").Value = "  -1.86%  "

$ws.Range("E23").Value = "  -2.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.91"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.88"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  -11.13%  "

$ws.Range("E28").Value = "  -5.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("E29").Value = "  +8.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -3.61%  "

$ws.Range("E31").Value = "  -4.37%  "

$ws.Range("E32").Value = "  -8.69%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.43"
$ws.Range("E35").Value = "  -3.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").Value = "  -3.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  -4.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.80"
$ws.Range("E38").Value = "  -1.56%  "

$ws.Range("E39").Value = "  -4.21%  "

$ws.Range("E40").Value = "  -6.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.15"
$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("E42").Value = "  -3.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.28"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("E44").Value = "  -4.04%  "

$ws.Range("D45").Value = "2.680.13"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.08"
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "342.67"
$ws.Range("E48").Value = "  -3.71%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  -1.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.49"
$ws.Range("E51").Value = "  -4.76%  "
